$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39:70 down to 40:71
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with its data (matches the row below it for the
# columns that stay the same, with the changed columns per the diff).
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44579
$ws.Range("D39").NumberFormat = $ws.Range("D40").NumberFormat
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103002
$ws.Range("J39").Value = "Ciruela"
$ws.Range("K39").Value = "Black Amber"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 190
$ws.Range("N39").Value = 10000
$ws.Range("O39").Value = 10000
$ws.Range("P39").Value = 10000
$ws.Range("Q39").Value = "$/bandeja 18 kilos granel"
$ws.Range("R39").Value = "Provincia de Curicó"
$ws.Range("S39").Value = 556
$ws.Range("T39").Value = 18
